# Adding instantiation of attack gameobjects
# - Inserts 4 new columns (D:G) on the Units sheet to make room for
#   individual attack slot columns (attack1..attack5, where the former
#   single "attackName" column becomes the first of five).
# - Relabels the header row and the default-unit data row accordingly.
# - Appends a new actMoveSpeed-style block (Max/Min/Locked) at the end of
#   the table (columns BT:BW) for the new 5th attack slot's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 4 blank columns before column D (former "actHealth" column),
# shifting the attributes table right. This turns the old single-column
# "attackName" field (C) into the start of a 5-column attack1..attack5
# block (C:G).
$ws.Columns("D:G").Insert()

# --- Row 1 (headers) ---------------------------------------------------
$ws.Range("C1").Value = "attack1"
$ws.Range("D1").Value = "attack2"
$ws.Range("E1").Value = "attack3"
$ws.Range("F1").Value = "attack4"
$ws.Range("G1").Value = "attack5"

# --- Row 2 (Default Unit data) -----------------------------------------
# C2 previously held "defaultAttack" (shifted from the old index after the
# sst churn); keep that as the default value for the first attack slot.
# D2:G2 (attack2..attack5) are left blank for the default unit.
$ws.Range("C2").Value = "defaultAttack"

# --- New trailing actMoveSpeed block (columns BT:BW) --------------------
$ws.Range("BT1").Value = "actMoveSpeed"
$ws.Range("BU1").Value = "actMoveSpeedMax"
$ws.Range("BV1").Value = "actMoveSpeedMin"
$ws.Range("BW1").Value = "actMoveSpeedLocked"

$ws.Range("BT2").Value = 10
$ws.Range("BU2").Value = 999
$ws.Range("BV2").Value = 1
$ws.Range("BW2").Value = $false

# Move the active selection to F2, matching where the author was editing.
$ws.Range("F2").Select()
